$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Remove the "License Information" Heading 2 paragraph entirely.
# ---------------------------------------------------------------------
$rng = $d.Content
$null = $rng.Find.Execute("License Information", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng.Expand(4)
$rng.Delete()

# ---------------------------------------------------------------------
# 2. Remove the "This PDF version is provided under the same license."
#    paragraph entirely (it merges into the preceding paragraph).
# ---------------------------------------------------------------------
$rng2 = $d.Content
$null = $rng2.Find.Execute("This PDF version is provided under the same license.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng2.Expand(4)
$rng2.Delete()

# ---------------------------------------------------------------------
# 3. Rename the bold "Termos Chave (Biblica)" run (the one that starts
#    the license paragraph, identified via its unique trailing text) to
#    "Biblica Study Notes (Key Terms)", preserving its bold formatting.
# ---------------------------------------------------------------------
$anchor = $d.Content
$null = $anchor.Find.Execute("Termos Chave (Biblica) (Portuguese)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boldRng = $d.Range($anchor.Start, $anchor.Start)
$null = $boldRng.Find.Execute("Termos Chave (Biblica)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$boldRng.Text = "Biblica Study Notes (Key Terms)"
$boldEnd = $boldRng.End

# ---------------------------------------------------------------------
# 4. Replace the remainder of that paragraph (the old "is based on: ..."
#    sentence plus the two hyperlinks) with the new license / adaptation
#    text, as plain (non-hyperlinked) runs.
# ---------------------------------------------------------------------
$p = $d.Paragraphs(4)
$pEnd = $p.Range.End
$tailRng = $d.Range($boldEnd, $pEnd - 1)
$tailRng.Text = " © 2023 Biblica Inc. Released under CC BY-SA 4.0 license. Biblica Study Notes has been adapted in the following languages: Tok Pisin, Arabic (عربي), French (Français), Hindi (हिंदी), Indonesian (Bahasa Indonesia), Portuguese (Português), Russian (Русский), Spanish (Español), Swahili (Kiswahili), and Simplified Chinese (简体中文)from Biblica Study Notes © 2023 Biblica Inc. Released under CC BY-SA 4.0 license by Mission Mutual."
$tailRng.Font.Bold = $false

# ---------------------------------------------------------------------
# 5. Remove the italic "Pai, Parábolas, Pastor, ..." key-terms listing
#    paragraph entirely.
# ---------------------------------------------------------------------
$rng3 = $d.Content
$null = $rng3.Find.Execute("Pai, Parábolas, Pastor, Patmos, Patriarca", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$rng3.Expand(4)
$rng3.Delete()
